$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The question/answer table previously started on row 2, leaving row 1 blank
# (dimension was A2:C102). Re-running the evaluation/export rewrote the sheet
# starting at row 1, so delete the leading blank row and let every row of
# data shift up by one (old row 2 -> row 1, ... old row 102 -> row 101).
$ws.Rows("1").Delete()

# The active cell/selection moved from B6 to A6 after the shift.
$ws.Range("A6").Select()
